# Rename worksheet from "Code Metrics" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Update Code Metrics cell values to reflect refactor (coverage improvements, legacy code cleanup)
$ws.Range("G2").Value = 344
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 1679
$ws.Range("K2").Value = 240
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 388
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 315
$ws.Range("G22").Value = 189
$ws.Range("I22").Value = 37
$ws.Range("J22").Value = 683
$ws.Range("K22").Value = 111
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 153
$ws.Range("J26").Value = 136
$ws.Range("J27").Value = 68
$ws.Range("J30").Value = 54
$ws.Range("I31").Value = 7
$ws.Range("J31").Value = 77
$ws.Range("J32").Value = 33
$ws.Range("J33").Value = 35
$ws.Range("F34").Value = 55
$ws.Range("G34").Value = 60
$ws.Range("I34").Value = 3
$ws.Range("K34").Value = 18
$ws.Range("F35").Value = 51
$ws.Range("G35").Value = 33
$ws.Range("I35").Value = 1
$ws.Range("K35").Value = 13
$ws.Range("J37").Value = 72
$ws.Range("J44").Value = 38
$ws.Range("F63").Value = 73
$ws.Range("I63").Value = 10
$ws.Range("J65").Value = 30
$ws.Range("J67").Value = 5
$ws.Range("J81").Value = 102
$ws.Range("J90").Value = 17
$ws.Range("J91").Value = 5
$ws.Range("G108").Value = 45
$ws.Range("J108").Value = 333
$ws.Range("K108").Value = 48
$ws.Range("I109").Value = 4
$ws.Range("J109").Value = 19
$ws.Range("F112").Value = 80
$ws.Range("G112").Value = 40
$ws.Range("J112").Value = 299
$ws.Range("K112").Value = 45
$ws.Range("F119").Value = 93
$ws.Range("G119").Value = 2
$ws.Range("I119").Value = 2
$ws.Range("J119").Value = 27
$ws.Range("K119").Value = 2
$ws.Range("F121").Value = 87
$ws.Range("G121").Value = 1
$ws.Range("I121").Value = 2
$ws.Range("J121").Value = 20
$ws.Range("K121").Value = 1
$ws.Range("F137").Value = 67
$ws.Range("G137").Value = 4
$ws.Range("J137").Value = 18
$ws.Range("K137").Value = 6
$ws.Range("I139").Value = 7
$ws.Range("J139").Value = 10
$ws.Range("I144").Value = 27
$ws.Range("J144").Value = 139
$ws.Range("J145").Value = 74
$ws.Range("J148").Value = 30
$ws.Range("I149").Value = 7
$ws.Range("J149").Value = 11
$ws.Range("I151").Value = 6
$ws.Range("J151").Value = 14
